$d = $word.ActiveDocument

# Locate the "Font-Size Property:" heading paragraph. It currently sits in
# its own paragraph immediately after a standalone page-break paragraph.
# We merge the page-break paragraph into the heading paragraph (deleting
# the paragraph mark between them) and then drop the inherited tab stop so
# the merged paragraph keeps the simple (no-tabs) formatting that the
# page-break paragraph had.

$c = $d.Content
$found = $c.Find.Execute("Font-Size Property:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingStart = $c.Start

$paras = $d.Paragraphs
$headingIndex = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Start -eq $headingStart) {
        $headingIndex = $i
    }
}

$prevPara = $paras.Item($headingIndex - 1)
$mergeStart = $prevPara.Range.End - 1
$mergeEnd = $prevPara.Range.End

$markRange = $d.Range($mergeStart, $mergeEnd)
$markRange.Delete()

$paras2 = $d.Paragraphs
$mergedPara = $paras2.Item($headingIndex - 1)
$mergedPara.TabStops.ClearAll()
